$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$sortRange = $ws.Range("A2:C13")
$keyRange = $ws.Range("A2:A13")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 1, 0, 0) | Out-Null

$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.SortMethod = 1
$ws.Sort.Apply()

$ws.Range("B8").Select() | Out-Null
